$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell "time_taken" in F1, copying the header style from E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

# Populate the time_taken values for each data row (F2:F47)
$timeTaken = @{
    2 = "2021-10-05 13:41:04.606736"
    3 = "2021-10-05 13:41:04.606748"
    4 = "2021-10-05 13:41:04.606752"
    5 = "2021-10-05 13:41:04.606756"
    6 = "2021-10-05 13:41:04.606759"
    7 = "2021-10-05 13:41:04.606763"
    8 = "2021-10-05 13:41:04.606766"
    9 = "2021-10-05 13:41:04.606769"
    10 = "2021-10-05 13:41:04.606772"
    11 = "2021-10-05 13:41:04.606775"
    12 = "2021-10-05 13:41:04.606779"
    13 = "2021-10-05 13:41:04.606782"
    14 = "2021-10-05 13:41:04.606785"
    15 = "2021-10-05 13:41:04.606788"
    16 = "2021-10-05 13:41:04.606791"
    17 = "2021-10-05 13:41:04.606794"
    18 = "2021-10-05 13:41:04.606797"
    19 = "2021-10-05 13:41:04.606800"
    20 = "2021-10-05 13:41:04.606804"
    21 = "2021-10-05 13:41:04.606807"
    22 = "2021-10-05 13:41:04.606810"
    23 = "2021-10-05 13:41:04.606813"
    24 = "2021-10-05 13:41:04.606816"
    25 = "2021-10-05 13:41:04.606819"
    26 = "2021-10-05 13:41:04.606823"
    27 = "2021-10-05 13:41:04.606826"
    28 = "2021-10-05 13:41:04.606829"
    29 = "2021-10-05 13:41:04.606832"
    30 = "2021-10-05 13:41:04.606835"
    31 = "2021-10-05 13:41:04.606839"
    32 = "2021-10-05 13:41:04.606842"
    33 = "2021-10-05 13:41:04.606845"
    34 = "2021-10-05 13:41:04.606849"
    35 = "2021-10-05 13:41:04.606852"
    36 = "2021-10-05 13:41:04.606855"
    37 = "2021-10-05 13:41:04.606858"
    38 = "2021-10-05 13:41:04.606861"
    39 = "2021-10-05 13:41:04.606864"
    40 = "2021-10-05 13:41:04.606867"
    41 = "2021-10-05 13:41:04.606870"
    42 = "2021-10-05 13:41:04.606874"
    43 = "2021-10-05 13:41:04.606877"
    44 = "2021-10-05 13:41:04.606880"
    45 = "2021-10-05 13:41:04.606883"
    46 = "2021-10-05 13:41:04.606886"
    47 = "2021-10-05 13:41:04.606889"
}

foreach ($row in $timeTaken.Keys) {
    $ws.Cells.Item($row, 6).Value = $timeTaken[$row]
}

Write-Output "time_taken column added"
